$wb = $excel.ActiveWorkbook

# ALC row 3: One for the Books
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 40000
$ws.Range("J3").Value = 40000
$ws.Range("L3").Value = 40000
$ws.Range("N3").Value = -40228

# ALC row 40: Stuck in the Moment
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2229.077
$ws.Range("I40").Value = 2297.9
$ws.Range("J40").Value = 1999.6666
$ws.Range("K40").Value = 2297.9
$ws.Range("L40").Value = 1999.6666
$ws.Range("M40").Value = -2122.9
$ws.Range("N40").Value = -2349.6666

# ALC row 100: Asking for a Friend
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 111114450
$ws.Range("I100").Value = 333333340
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 333333340
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -333332799
$ws.Range("N100").Value = -6082

# ALC row 102: Spell-rebound
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H102").Value = 40000
$ws.Range("J102").Value = 40000
$ws.Range("L102").Value = 40000
$ws.Range("N102").Value = -46490

# ARM row 2: Ain't Got No Ingots
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2710.625
$ws.Range("I2").Value = 2730.8333
$ws.Range("K2").Value = 2730.8333
$ws.Range("M2").Value = -2617.8333

# ARM row 31: I Was a Teenage Wailer
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 15065.556
$ws.Range("I31").Value = 11866.667
$ws.Range("J31").Value = 21463.334
$ws.Range("K31").Value = 11866.667
$ws.Range("L31").Value = 21463.334
$ws.Range("M31").Value = -11572.667
$ws.Range("N31").Value = -22051.334

# ARM row 116: No Scope
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2710.625
$ws.Range("I116").Value = 2730.8333
$ws.Range("K116").Value = 2730.8333
$ws.Range("M116").Value = -436.8332999999998

# ARM row 132: Don't Bore Me, Ore Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2153
$ws.Range("I132").Value = 782.25806
$ws.Range("J132").Value = 4808.8125
$ws.Range("K132").Value = 2346.77418
$ws.Range("L132").Value = 14426.4375
$ws.Range("M132").Value = 183.2258200000001
$ws.Range("N132").Value = -19486.4375

# BSM row 3: Hells Bells
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2710.625
$ws.Range("I3").Value = 2730.8333
$ws.Range("K3").Value = 2730.8333
$ws.Range("M3").Value = -2616.8333

# BSM row 86: Through Thick and Thin
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3221
$ws.Range("I86").Value = 3776.25
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 3776.25
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = -2653.25
$ws.Range("N86").Value = -3246

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3221
$ws.Range("I89").Value = 3776.25
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 18881.25
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = -13265.25
$ws.Range("N89").Value = -16232

# BSM row 94: High Steal
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 825.4783
$ws.Range("I94").Value = 912
$ws.Range("J94").Value = 713
$ws.Range("K94").Value = 912
$ws.Range("L94").Value = 713
$ws.Range("M94").Value = -461
$ws.Range("N94").Value = -1615

# BSM row 99: Meddle in Metal
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1386.2858
$ws.Range("I99").Value = 907.5
$ws.Range("J99").Value = 2024.6666
$ws.Range("K99").Value = 907.5
$ws.Range("L99").Value = 2024.6666
$ws.Range("M99").Value = 590.5
$ws.Range("N99").Value = -5020.6666

# BSM row 105: Ingot to Wing It
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2438.9524
$ws.Range("I105").Value = 2060.9
$ws.Range("K105").Value = 2060.9
$ws.Range("M105").Value = -313.9000000000001

# CRP row 16: Raise the Roof
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1448.2778
$ws.Range("I16").Value = 1603
$ws.Range("J16").Value = 1046
$ws.Range("K16").Value = 1603
$ws.Range("L16").Value = 1046
$ws.Range("M16").Value = -1316
$ws.Range("N16").Value = -1620

# CRP row 113: Patient Patients
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1448.2778
$ws.Range("I113").Value = 1603
$ws.Range("J113").Value = 1046
$ws.Range("K113").Value = 1603
$ws.Range("L113").Value = 1046
$ws.Range("M113").Value = 567
$ws.Range("N113").Value = -5386

# GSM row 70: Sky Is the Limit
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 39665.855
$ws.Range("I70").Value = 127504
$ws.Range("J70").Value = 4530.6
$ws.Range("K70").Value = 127504
$ws.Range("L70").Value = 4530.6
$ws.Range("M70").Value = -127234
$ws.Range("N70").Value = -5070.6

# GSM row 73: Hulls of Broken Dreams (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 39665.855
$ws.Range("I73").Value = 127504
$ws.Range("J73").Value = 4530.6
$ws.Range("K73").Value = 127504
$ws.Range("L73").Value = 4530.6
$ws.Range("M73").Value = -126568
$ws.Range("N73").Value = -6402.6

# LTW row 40: Best Served Toad
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 65505
$ws.Range("I40").Value = 101438
$ws.Range("J40").Value = 5616.6665
$ws.Range("K40").Value = 101438
$ws.Range("L40").Value = 5616.6665
$ws.Range("M40").Value = -101302
$ws.Range("N40").Value = -5888.6665

# LTW row 42: Slave to Fashion
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 6666.3335
$ws.Range("J42").Value = 6666.3335
$ws.Range("L42").Value = 6666.3335
$ws.Range("N42").Value = -7792.3335

# LTW row 46: Supply Side Logic
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 78846.16
$ws.Range("J46").Value = 1349.9
$ws.Range("L46").Value = 1349.9
$ws.Range("N46").Value = -1725.9

# LTW row 49: First They Came for the Heretics
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H49").Value = 6666.3335
$ws.Range("J49").Value = 6666.3335
$ws.Range("L49").Value = 6666.3335
$ws.Range("N49").Value = -6960.3335

# LTW row 68: You Could Say It's a Moving Target
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1878.75
$ws.Range("I68").Value = 1530
$ws.Range("J68").Value = 2088
$ws.Range("K68").Value = 1530
$ws.Range("L68").Value = 2088
$ws.Range("M68").Value = -781
$ws.Range("N68").Value = -3586

# LTW row 71: They Call It Bloody Mary (L)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1878.75
$ws.Range("I71").Value = 1530
$ws.Range("J71").Value = 2088
$ws.Range("K71").Value = 7650
$ws.Range("L71").Value = 10440
$ws.Range("M71").Value = -3906
$ws.Range("N71").Value = -17928

# LTW row 82: Trainin' the Neck
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2884.6155
$ws.Range("I82").Value = 1225
$ws.Range("J82").Value = 3622.2222
$ws.Range("K82").Value = 1225
$ws.Range("L82").Value = 3622.2222
$ws.Range("M82").Value = -864
$ws.Range("N82").Value = -4344.2222

# LTW row 85: Training Is Only Skintight (L)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2884.6155
$ws.Range("I85").Value = 1225
$ws.Range("J85").Value = 3622.2222
$ws.Range("K85").Value = 1225
$ws.Range("L85").Value = 3622.2222
$ws.Range("M85").Value = 23
$ws.Range("N85").Value = -6118.2222

# LTW row 132: Tenets of Tanning
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 14969.05
$ws.Range("I132").Value = 18952
$ws.Range("K132").Value = 56856
$ws.Range("M132").Value = -54326

# WVR row 113: A Tender Table
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1020.3077
$ws.Range("I113").Value = 939.7368
$ws.Range("K113").Value = 2819.2104
$ws.Range("M113").Value = -649.2103999999999
